$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '67.125.91'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '3.131.48'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '581.17'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '174.36'
$ws.Range('E6').Value = '  +1.07%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.522'
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.155'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').Value = '6.43'
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').Value = '0.482'
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').Value = '0.0000250'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '37.67'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '0.122'
$ws.Range('E14').Value = '  -1.60%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.652.91'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '67.109.42'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '7.17'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.130.94'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '16.40'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '492.63'
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').Value = '0.712'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '7.91'
$ws.Range('E22').Value = '  +4.57%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '84.31'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '13.38'
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = '2.30'
$ws.Range('E25').Value = '  -2.78%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').Value = '10.46'
$ws.Range('E26').Value = '  +3.83%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('B28').Value = 'NEARProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D28').Value = '7.97'
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '2.36'
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '2.70'
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '28.73'
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.115'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '0.0₃0952'
$ws.Range('E33').Value = '  -6.17%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = '5.94'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').Value = '0.979'
$ws.Range('E36').Value = '  -2.84%  '
$ws.Range('B37').Value = 'Arweave'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D37').Value = '46.61'
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '2.06'
$ws.Range('E38').Value = '  -3.03%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '50.16'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '0.313'
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.124'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Value = '8.58'
$ws.Range('E42').Value = '  -1.26%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.841.94'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = '387.09'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '2.60'
$ws.Range('E45').Value = '  -6.91%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0355'
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '135.90'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '25.07'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '2.23'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.108'
$ws.Range('E51').Value = '  -0.28%  '
